# Seguimiento Pruebas CU.docx
# Row "Generar Programa PDF" -> mark Regression test as performed/approved:
#  - shade every cell of the row green (00B050) instead of red (D99594)
#  - update several of the cell texts to reflect completion

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the row whose first cell reads "Generar Programa PDF" (row 13
# in this table), rather than hard-coding the index, to stay robust.
$targetRow = $null
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $firstCellText = $t.Rows.Item($i).Cells.Item(1).Range.Text
    if ($firstCellText -like "Generar Programa PDF*") {
        $targetRow = $t.Rows.Item($i)
        break
    }
}

# 00B050 as a Word BGR-packed long (R + G*256 + B*65536)
$green = 0x00 + (0xB0 * 256) + (0x50 * 65536)

# New text per cell index (1-based); $null means "leave text unchanged"
$newText = @{
    3  = "Realizada"
    5  = "18 y 19/06/2020"
    7  = "N/A"
    8  = "SI"
    9  = "N/A"
    10 = "NO"
    11 = "Aprobado (CU cerrado - No se debe modificar)"
}

for ($i = 1; $i -le $targetRow.Cells.Count; $i++) {
    $cell = $targetRow.Cells.Item($i)
    $cell.Shading.BackgroundPatternColor = $green
    if ($newText.ContainsKey($i)) {
        $cell.Range.Text = $newText[$i]
    }
}
